$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Comment.Delete()
$ws.Range("D1").ClearFormats()
$ws.Range("D1:F2").Value = "asdfasdf"
$ws.Range("D1:E2").Select()
